$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I (shifts IC4..3MPent from I..O to J..P)
$ws.Columns("I:I").Insert()

# New column header: C16H34
$ws.Range("I1").Value = "C16H34"

# Correct a data value that was wrong before the shift (row 2, was RUNNING -> now YES)
# After the column insert this cell lives at P2 (was O2).
$ws.Range("P2").Value = "YES"

# Re-select the cell the author left active
$ws.Range("P6").Select()
